$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark that sits right after "MP73010" in the
#    title line (it moves further down the document in the edited version).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Replace the (until now empty) paragraph right before the very last
#    paragraph with the new "version management" blurb, broken up into the
#    same run/bookmark structure produced by Word when the text was typed
#    interactively (several small runs, with a fresh "_GoBack" bookmark
#    landing after "management ").
$targetParagraph = $d.Paragraphs(6)
$insertionPoint = $d.Range($targetParagraph.Range.Start, $targetParagraph.Range.Start)

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:t xml:space="preserve">The term </w:t></w:r>' +
    '<w:r><w:t>v</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">ersion </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">management </w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>refer to the means of</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> effectively track</w:t></w:r>' +
    '<w:r><w:t>ing</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and control</w:t></w:r>' +
    '<w:r><w:t>ling</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> changes</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> to a </w:t></w:r>' +
    '<w:r><w:t>group of entities usually files and the information contained within them</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
    '</w:p>'

$insertionPoint.InsertXML($newParagraphXml)
